$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("event")

# Row 4 updates
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("P4").Value = 0

# Row 10 updates
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0

# Update the active selection to T10 on the event sheet
$ws.Activate()
$ws.Range("T10").Select()
